# Code Merge Changes - 9/25/2017
# Update the "TestResultExcelFilePath" value (cell H2) on each SSP/report
# worksheet to point at the new automation results location.

$wb = $excel.ActiveWorkbook

$newPath = "F:\\Automation_TestResults\\Payroll_Tax_StatutoryScenarios\\Automation Test Result for Statutory Scenarios201718.xlsx"

# First sheet keeps a trailing newline in the cell content (as authored),
# the rest use the plain value.
$ws = $wb.Worksheets.Item("ProcessPayrolFor20FourWeeklySSP")
$ws.Range("H2").Value = $newPath + "`n"

$sheetNames = @(
    "ProcessPayrolFor24FourWeeklySSP",
    "ProcessPayrolFor28FourWeeklySSP",
    "AverageWeeklyEarningsTestReport",
    "ProcessPayrolFor32FourWeeklySSP",
    "ProcessPayrolFor36FourWeeklySSP",
    "AverageWeeklyEarningsTestReprt2",
    "ProcessPayrolFor40FourWeeklySSP"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("H2").Value = $newPath
}
